# Weekly update: a new daily price record for Ajo (Chino, Primera) at
# Terminal Hortofrutícola Agro Chillán is inserted as row 63, pushing the
# existing rows 63-182 down to 64-183 (dimension grows from A1:R182 to
# A1:R183).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 63 (shifts rows 63.. down by one).
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(63, 1).Value = 7
$ws.Cells.Item(63, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(63, 3).Value = "Ñuble"
$ws.Cells.Item(63, 4).Value = 44581
$ws.Cells.Item(63, 5).Value = 16
$ws.Cells.Item(63, 6).Value = 100112003
$ws.Cells.Item(63, 7).Value = "Ajo"
$ws.Cells.Item(63, 8).Value = "Chino"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 60
$ws.Cells.Item(63, 11).Value = 20000
$ws.Cells.Item(63, 12).Value = 21000
$ws.Cells.Item(63, 13).Value = 20500
$ws.Cells.Item(63, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(63, 15).Value = "China"
$ws.Cells.Item(63, 16).Value = 2050
$ws.Cells.Item(63, 17).Value = 10
$ws.Cells.Item(63, 18).Value = "Hortaliza"

# Keep the date cell formatted like the rest of column D (yyyy-mm-dd date
# style, same as the cell it displaced).
$ws.Cells.Item(63, 4).NumberFormat = $ws.Cells.Item(64, 4).NumberFormat
